# Daily attendance processing - 2025-10-30 04:47:56
# Reorders the "Recorded By" (column G) entries for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Literal mapping of old -> new "Recorded By" text, as produced by the
# upstream processing job (reorders the names/emails within the cell).
$map = @{}
$map["System, dnasr281@gmail.com"] = "dnasr281@gmail.com, System"
$map["backup@backdoor.com, System"] = "System, backup@backdoor.com"
$map["backup@backdoor.com, System, system"] = "System, backup@backdoor.com, system"
$map["admin@admin.com, dnasr281@gmail.com"] = "dnasr281@gmail.com, admin@admin.com"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
